# Fruta / hortaliza, semanal
# Insert a new weekly record at row 21, shifting existing rows 21-49 down to 22-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21; this shifts rows 21-49 down to 22-50,
# carrying their values and formatting (e.g. the date style in column D) with them.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with this week's data.
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C21").Value = "Arica y Parinacota"
$ws.Range("D21").Value = 44987
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100104
$ws.Range("H21").Value = "Frutos de pepita"
$ws.Range("I21").Value = 100104005
$ws.Range("J21").Value = "Pera"
$ws.Range("K21").Value = "Packham's Triumph"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 19000
$ws.Range("P21").Value = 18500
$ws.Range("Q21").Value = "`$/caja 20 kilos granel"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 925
$ws.Range("T21").Value = 20
